$wb = $excel.ActiveWorkbook
$n = $wb.Names.Item("_xlchart.v2.0")
$n.RefersTo = "=Sheet2!`$E`$32:`$E`$35"
Write-Host "new refersto:" $n.RefersTo
